$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-order / relabel header row:
# A1 = IP, B1 = Usuario, C1 = Password, D1 = Puerto API, E1 = Nombre del Equipo
$ws.Range("A1").Value = "IP"
$ws.Range("B1").Value = "Usuario"
$ws.Range("C1").Value = "Password"
$ws.Range("D1").Value = "Puerto API"
$ws.Range("E1").Value = "Nombre del Equipo"

# Update the API port values in column D
$ws.Range("D2").Value = 4444
$ws.Range("D3").Value = 4444

# Move the active selection to D3 (matches the saved selection in the sheet view)
$ws.Range("D3").Select()
